$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("short term")

# New TODO items to append at the bottom of the "short term" list
$task83 = "83. try to make the size scale of nb_of_items print only integer values in sensitivityplot and trendplot"
$task84 = "84. ajouter dans le package une fonction de visualisation de l'ensemble des données en ACP pour voir la cohérence et détecter d'éventuels outliers"

# Row 29: task 83, assigned to "A" -- copy style from row 26 (A: fillId orange, B: plain)
$ws.Range("A26:B26").Copy() | Out-Null
$ws.Range("A29:B29").PasteSpecial(-4122) | Out-Null

$ws.Range("A29").Value = $task83
$ws.Range("B29").Value = "A"

# Row 30: task 84, assigned to "ML" -- copy style from row 25 (A: fillId yellow, B: plain)
$ws.Range("A25:B25").Copy() | Out-Null
$ws.Range("A30:B30").PasteSpecial(-4122) | Out-Null

$ws.Range("A30").Value = $task84
$ws.Range("B30").Value = "ML"

$excel.CutCopyMode = 0

# Update view so the newly added row is visible/selected, mirroring the saved view state
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("A30").Select() | Out-Null

$wb.Save()
